$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the capitalisation of the "depth" header in K1 -> "Depth"
$ws.Range("K1").Value = "Depth"

# Column K used to be a shared formula (A*-5); replace it with the new
# literal depth values (the underlying relationship is now 4*A-11, but it
# is stored as plain numbers rather than a formula).
$kValues = @(-7, -7, -3, -3, 1, 1, 5, 5, 9, 9, 13, 13, 17, 17, 21, 21, 25, 25, 29, 29, 33, 33, 37, 37, 41, 41, 45, 45)
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 11).Value = $kValues[$i]
}

# Move/update the active selection to where editing left off.
$ws.Range("I35").Select()
